$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Industries") values for rows 18-176 were changed from 1 to 0
$ws.Range("H18:H176").Value = 0
